$d = $word.ActiveDocument

# All four "Dates de la campanya ..." paragraphs (one of which has a stray
# leading space in its own run) get collapsed into a single plain run with
# the translated/updated campaign dates text and no run formatting.
$searchPattern = "[ ]{0,1}Dates de la campanya*novembre de desembre 8"
$newText = "Dates de la campanya Perseus: 16-25 de gener, del 7 al 16 de novembre, del 6 al 15 de desembre"

$rng = $d.Content
$rng.Find.ClearFormatting()
$rng.Find.MatchWildcards = $true

$count = 0
while ($rng.Find.Execute($searchPattern, $false, $false, $true, $false, $false, $true, 1, $false, "", 0)) {
    $count = $count + 1

    # Wipe out the matched text, then insert the replacement as a brand new
    # run (InsertAfter on an emptied range drops any inherited rPr), so the
    # resulting run carries no formatting - matching the target markup.
    $rng.Text = ""
    $rng.InsertAfter($newText)
    $rng.Collapse(0)

    if ($count -ge 20) { break }
}
